$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the existing header formatting (bold font, thin border, centered)
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill data rows 2 through 11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20180335
    $ws.Cells.Item($r, 15).Value = 8
}
